# Transition rule: add "Within 5 miles" and "Within 10 miles" of HFC
# production facility columns (F and G) to both the "Means" and
# "Standard Deviations" sheets, and update a couple of recalculated
# summary values that shifted as a result of adding the new radii.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Means"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cells for the two additional distance columns.
$ws1.Cells.Item(1, 6).Value2 = "Within 5 miles of HFC production facility"
$ws1.Cells.Item(1, 7).Value2 = "Within 10 miles of HFC production facility"

# New data columns F (5 mi) and G (10 mi) for rows 2-10 all hold #NUM!
# errors, matching the existing C/D/E columns.
for ($r = 2; $r -le 10; $r++) {
    $ws1.Cells.Item($r, 6).Value2 = "#NUM!"
    $ws1.Cells.Item($r, 7).Value2 = "#NUM!"
}

# Recalculated summary values that changed with the new radii included.
$ws1.Cells.Item(9, 2).Value2 = 29      # Total Cancer Risk (per million): 32 -> 29
$ws1.Cells.Item(10, 2).Value2 = 0.37   # Total Respiratory (hazard quotient): 0.44 -> 0.37

# ----------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New header cells for the two additional distance columns (SD variants).
$ws2.Cells.Item(1, 6).Value2 = "Within 5 mile of HFC production facility SD"
$ws2.Cells.Item(1, 7).Value2 = "Within 10 mile of HFC production facility SD"

# New data columns F (5 mi) and G (10 mi) for rows 2-10 all hold 0,
# matching the existing C/D/E columns.
for ($r = 2; $r -le 10; $r++) {
    $ws2.Cells.Item($r, 6).Value2 = 0
    $ws2.Cells.Item($r, 7).Value2 = 0
}

# Recalculated summary values that changed with the new radii included.
$ws2.Cells.Item(9, 2).Value2 = 10      # Total Cancer Risk (per million) SD: 12 -> 10
$ws2.Cells.Item(10, 2).Value2 = 0.14   # Total Respiratory (hazard quotient) SD: 0.15 -> 0.14
